$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new data rows after the existing data block (row 12),
# shifting the old Subtotal/Shipping/Total rows (14,15,16) down to (17,18,19)
# and leaving row 16 blank, matching the target layout.
$ws.Rows("13:15").Insert()

# --- Row 13: Poster Print ---
$ws.Range("B13").Value = "Poster Print"
$ws.Range("C13").Value = "Final Poster Print"
$ws.Range("D13").Value = "ENGR Computer Lab"
$ws.Range("E13").Formula = "=17.94+2.03"
$ws.Range("F13").Value = 1
$ws.Range("G13").Formula = "=F13*E13"
$ws.Range("H13").Value = 0

# --- Row 14: Transistor ---
$ws.Range("B14").Value = "Transistor"
$ws.Range("D14").Value = "ECE Store"
$ws.Range("C14").Value = "Transistor laser driver"
$ws.Range("E14").Value = 4.6399999999999997
$ws.Range("F14").Value = 1
$ws.Range("G14").Formula = "=F14*E14"
$ws.Range("H14").Value = 0

# --- Row 15: Flexmod P3 (backup laser driver) ---
$ws.Range("B15").Value = "Flexmod P3"
$ws.Range("D15").Value = "Innolasers"

# --- Fix row 7 description (was "Laser Driver Heat Sink") ---
$ws.Range("C7").Value = "Laser Driver"

$ws.Range("C15").Value = "Backup Laser Driver"
$ws.Range("E15").Value = 35.99
$ws.Range("F15").Value = 1
$ws.Range("G15").Formula = "=F15*E15"
$ws.Range("H15").Value = 35.299999999999997

# Extend the shared formula's range comment on G4 down through the new rows.
$ws.Range("G4:G15").FormulaR1C1 = "=RC6*RC5"

# --- Subtotal / Shipping / Total rows now live at 17-19 ---
$ws.Range("G17").Formula = "=SUM(G3:G15)"
$ws.Range("G18").Formula = "=SUM(H3:H15)"
$ws.Range("G19").Formula = "=SUM(G17:H18)"

# Column D needs to widen to fit "ENGR Computer Lab" / "Transistor laser driver".
$ws.Columns("D:D").ColumnWidth = 18

# Sheet view: zoom to 85% and move the active selection to D20.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("D20").Select()
